# Updated mapping to prototype cart
$wb = $excel.ActiveWorkbook

# Update the device_cart mapping: "zovag" -> "getit"
$ws1 = $wb.Worksheets.Item("device_cart")
$ws1.Range("A2").Value = "getit"

# Make device_cart the active/selected sheet (was rfid_item)
$ws1.Activate()
